$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.38261866569519
$ws.Range("B1").Value = -1
$ws.Range("D1").Value = 0.5292857885360718
$ws.Range("E1").Value = 0.6987159252166748
